$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 16676487
$ws.Range("I62").Value = 25012280
$ws.Range("J62").Value = 4901.2
$ws.Range("K62").Value = 25012280
$ws.Range("L62").Value = 4901.2
$ws.Range("M62").Value = -25011656
$ws.Range("N62").Value = -6149.2

$ws.Range("H64").Value = 2881.818
$ws.Range("I64").Value = 3800
$ws.Range("J64").Value = 2537.5
$ws.Range("K64").Value = 3800
$ws.Range("L64").Value = 2537.5
$ws.Range("M64").Value = -3552
$ws.Range("N64").Value = -3033.5

$ws.Range("H65").Value = 16676487
$ws.Range("I65").Value = 25012280
$ws.Range("J65").Value = 4901.2
$ws.Range("K65").Value = 125061400
$ws.Range("L65").Value = 24506
$ws.Range("M65").Value = -125058280
$ws.Range("N65").Value = -30746

$ws.Range("H67").Value = 2881.818
$ws.Range("I67").Value = 3800
$ws.Range("J67").Value = 2537.5
$ws.Range("K67").Value = 3800
$ws.Range("L67").Value = 2537.5
$ws.Range("M67").Value = -2942
$ws.Range("N67").Value = -4253.5

$ws.Range("H98").Value = 125000910
$ws.Range("I98").Value = 208333860
$ws.Range("J98").Value = 1499.75
$ws.Range("K98").Value = 208333860
$ws.Range("L98").Value = 1499.75
$ws.Range("M98").Value = -208332362
$ws.Range("N98").Value = -4495.75

$ws.Range("H99").Value = 259.38235
$ws.Range("I99").Value = 244.96875
$ws.Range("K99").Value = 734.90625
$ws.Range("M99").Value = 763.09375

$ws.Range("H122").Value = 125000910
$ws.Range("I122").Value = 208333860
$ws.Range("J122").Value = 1499.75
$ws.Range("K122").Value = 625001580
$ws.Range("L122").Value = 4499.25
$ws.Range("M122").Value = -624999130
$ws.Range("N122").Value = -9399.25

$ws.Range("H132").Value = 8776789
$ws.Range("I132").Value = 10420602
$ws.Range("J132").Value = 9783.888999999999
$ws.Range("K132").Value = 31261806
$ws.Range("L132").Value = 29351.667
$ws.Range("M132").Value = -31259276
$ws.Range("N132").Value = -34411.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 778.9394
$ws.Range("I2").Value = 620.4737
$ws.Range("J2").Value = 994
$ws.Range("K2").Value = 620.4737
$ws.Range("L2").Value = 994
$ws.Range("M2").Value = -507.4737
$ws.Range("N2").Value = -1220

$ws.Range("H32").Value = 11239332
$ws.Range("I32").Value = 3280.2666
$ws.Range("K32").Value = 3280.2666
$ws.Range("M32").Value = -2993.2666

$ws.Range("H45").Value = 92135
$ws.Range("I45").Value = 143553.14
$ws.Range("J45").Value = 2153.25
$ws.Range("K45").Value = 143553.14
$ws.Range("L45").Value = 2153.25
$ws.Range("M45").Value = -143176.14
$ws.Range("N45").Value = -2907.25

$ws.Range("H61").Value = 1559.0454
$ws.Range("I61").Value = 1599.95
$ws.Range("J61").Value = 1150
$ws.Range("K61").Value = 1599.95
$ws.Range("L61").Value = 1150
$ws.Range("M61").Value = -1387.95
$ws.Range("N61").Value = -1574

$ws.Range("H97").Value = 1453.762
$ws.Range("I97").Value = 1407.6154
$ws.Range("J97").Value = 1528.75
$ws.Range("K97").Value = 1407.6154
$ws.Range("L97").Value = 1528.75
$ws.Range("M97").Value = -911.6153999999999
$ws.Range("N97").Value = -2520.75

$ws.Range("H102").Value = 1111
$ws.Range("I102").Value = 1178.8889
$ws.Range("J102").Value = 500
$ws.Range("K102").Value = 1178.8889
$ws.Range("L102").Value = 500
$ws.Range("M102").Value = 443.1111000000001
$ws.Range("N102").Value = -3744

$ws.Range("H110").Value = 1768.3334
$ws.Range("I110").Value = 830.7143
$ws.Range("K110").Value = 830.7143
$ws.Range("M110").Value = 1214.2857

$ws.Range("H116").Value = 778.9394
$ws.Range("I116").Value = 620.4737
$ws.Range("J116").Value = 994
$ws.Range("K116").Value = 620.4737
$ws.Range("L116").Value = 994
$ws.Range("M116").Value = 1673.5263
$ws.Range("N116").Value = -5582

$ws.Range("H132").Value = 817791.2
$ws.Range("I132").Value = 580.4308
$ws.Range("J132").Value = 8406177
$ws.Range("K132").Value = 1741.2924
$ws.Range("L132").Value = 25218531
$ws.Range("M132").Value = 788.7076000000002
$ws.Range("N132").Value = -25223591

$ws.Range("H136").Value = 1559.0454
$ws.Range("I136").Value = 1599.95
$ws.Range("J136").Value = 1150
$ws.Range("K136").Value = 4799.85
$ws.Range("L136").Value = 3450
$ws.Range("M136").Value = -2249.85
$ws.Range("N136").Value = -8550

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 778.9394
$ws.Range("I3").Value = 620.4737
$ws.Range("J3").Value = 994
$ws.Range("K3").Value = 620.4737
$ws.Range("L3").Value = 994
$ws.Range("M3").Value = -506.4737
$ws.Range("N3").Value = -1222

$ws.Range("H107").Value = 166667500
$ws.Range("I107").Value = 250000750
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 250000750
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -249998830
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 20769.092
$ws.Range("I105").Value = 26875
$ws.Range("J105").Value = 4486.6665
$ws.Range("K105").Value = 26875
$ws.Range("L105").Value = 4486.6665
$ws.Range("M105").Value = -25128
$ws.Range("N105").Value = -7980.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 916.8099999999999
$ws.Range("I131").Value = 507.5
$ws.Range("J131").Value = 933.86456
$ws.Range("K131").Value = 1522.5
$ws.Range("L131").Value = 2801.59368
$ws.Range("M131").Value = 3517.5
$ws.Range("N131").Value = -12881.59368

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2239.6365
$ws.Range("I102").Value = 2292.889
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2292.889
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -670.8890000000001
$ws.Range("N102").Value = -5244

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 507299.38
$ws.Range("I22").Value = 974108
$ws.Range("J22").Value = 1590
$ws.Range("K22").Value = 974108
$ws.Range("L22").Value = 1590
$ws.Range("M22").Value = -973813
$ws.Range("N22").Value = -2180

$ws.Range("H27").Value = 507299.38
$ws.Range("I27").Value = 974108
$ws.Range("J27").Value = 1590
$ws.Range("K27").Value = 974108
$ws.Range("L27").Value = 1590
$ws.Range("M27").Value = -974001
$ws.Range("N27").Value = -1804

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 39650
$ws.Range("J46").Value = 39650
$ws.Range("L46").Value = 39650
$ws.Range("N46").Value = -40112

$ws.Range("H96").Value = 2087.6875
$ws.Range("I96").Value = 1533.3334
$ws.Range("J96").Value = 2420.3
$ws.Range("K96").Value = 1533.3334
$ws.Range("L96").Value = 2420.3
$ws.Range("M96").Value = -160.3334
$ws.Range("N96").Value = -5166.3

$ws.Range("H132").Value = 25480.291
$ws.Range("I132").Value = 33964.562
$ws.Range("J132").Value = 8511.75
$ws.Range("K132").Value = 101893.686
$ws.Range("L132").Value = 25535.25
$ws.Range("M132").Value = -99363.68599999999
$ws.Range("N132").Value = -30595.25

$ws.Range("H134").Value = 39650
$ws.Range("J134").Value = 39650
$ws.Range("L134").Value = 118950
$ws.Range("N134").Value = -124020

$ws.Range("H136").Value = 9260459
$ws.Range("I136").Value = 12821295
$ws.Range("J136").Value = 2286.3333
$ws.Range("K136").Value = 38463885
$ws.Range("L136").Value = 6858.999899999999
$ws.Range("M136").Value = -38461335
$ws.Range("N136").Value = -11958.9999
